# Edit: slide 123 ("Live as if..." quote slide)
#  - move the picture down slightly
#  - widen/reposition the quote textbox
#  - replace the two-line quote with a new single-line quote, and
#    re-color the text from maroon (B22251) to burlywood (DEB887)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(123)

# Shape 1: "Picture 2" - slide the picture down (Top only changes)
$pic = $s.Shapes.Item(1)
$pic.Top = 156

# Shape 2: "Rectangle 2" - the quote textbox
$rect = $s.Shapes.Item(2)
$rect.Top = 11.864094488188977
$rect.Width = 696

$tr = $rect.TextFrame.TextRange
$tr.Text = "`"If someone is strong enough to bring you down, show them you are strong enough to get up.`""
$tr.Font.Color.RGB = 8894686
$tr.LanguageID = "en-US"
